# Bug fix (column location):
# The "seg_B" column had been written into column F instead of column B,
# which pushed every other column (Unnamed: 1 / date / report / 備註) one
# slot to the left of where it belongs. Shift B:E back into C:F (their
# correct position) and restore the seg_B header in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the old B:E block (Unnamed: 1, date, report, 備註 + all data rows)
# one column to the right, into C:F, so it lines up under the correct
# headers again.
$ws.Range("B1:E8").Copy($ws.Range("C1"))

# Column B is now free for the real "seg_B" column: clear any leftover
# values (the old numeric duplicates) and write the header back.
$ws.Range("B1:B8").ClearContents()
$ws.Range("B1").Value = "seg_B"

# Column F held the old (duplicate) seg_B numbers in the data rows; now
# that row 1's "備註" has been copied into F1, just clear out the leftover
# duplicate values in the data rows F2:F8.
$ws.Range("F2:F8").ClearContents()
